# Fix elevation units: rows 5-10 of column D (ELEVATION (FT)) were entered
# in meters by mistake; convert them to feet using the standard
# meters -> feet factor (1 m = 3.28084 ft), matching the corrected values
# used for the already-converted rows above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$metersToFeet = 3.28084

for ($row = 5; $row -le 10; $row++) {
    $cell = $ws.Cells.Item($row, 4)
    $meters = $cell.Value2
    $cell.Value = $meters * $metersToFeet
}
